$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (51-55) to append below the existing data (A1:H50).
# Each row mirrors the layout of the existing rows: Date (A, styled like
# the rest of column A), city (B), type (C), size (D), price (E),
# unit_price (F), land_size (G, sometimes blank), count (H).
$rows = @(
    @{ Row=51; Date=43893; City="Algyo";      Type="House"; Size=113.8095238095238; Price=31133333.33333333; UnitPrice=287581.7671387243; LandSize=464.8095238095238; Count=21 }
    @{ Row=52; Date=43893; City="Morahalom";  Type="House"; Size=122.4047619047619; Price=30370238.0952381;  UnitPrice=251363.0346858636; LandSize=153.3333333333333; Count=42 }
    @{ Row=53; Date=43893; City="Szeged";     Type="Flat";  Size=65.29819277108433; Price=27861588.85542169; UnitPrice=431660.3488087555; LandSize=$null;             Count=2656 }
    @{ Row=54; Date=43893; City="Szeged";     Type="Garage";Size=17.68965517241379; Price=3981839.08045977;  UnitPrice=233766.3635384706; LandSize=$null;             Count=87 }
    @{ Row=55; Date=43893; City="Szeged";     Type="House"; Size=162.2271062271062; Price=56672832.72283272; UnitPrice=974787.8027311168; LandSize=466.4029304029304; Count=819 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the style (date format, font, border, alignment) of the cell
    # directly above in column A so the new date cell matches the rest
    # of the column without introducing any new style definitions.
    $ws.Range("A" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum).PasteSpecial(-4122)

    $ws.Range("A" + $rowNum).Value = $r.Date
    $ws.Range("B" + $rowNum).Value = $r.City
    $ws.Range("C" + $rowNum).Value = $r.Type
    $ws.Range("D" + $rowNum).Value = $r.Size
    $ws.Range("E" + $rowNum).Value = $r.Price
    $ws.Range("F" + $rowNum).Value = $r.UnitPrice
    if ($null -ne $r.LandSize) {
        $ws.Range("G" + $rowNum).Value = $r.LandSize
    }
    $ws.Range("H" + $rowNum).Value = $r.Count
}

$excel.CutCopyMode = $false
